$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying weekly records (rows 2-9) get reshuffled: each row keeps its
# "identity" columns (market, product, etc.) but the data columns
# (Fecha, Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Precio $/Kg) are redistributed across the rows as follows:
#   new row 2 <- old row 3
#   new row 3 <- old row 6
#   new row 4 <- old row 9
#   new row 5 <- old row 2
#   new row 6 <- old row 7
#   new row 7 <- old row 5
#   new row 8 <- old row 4
#   new row 9 <- old row 8

# Capture current (old) values for the affected columns before overwriting.
$oldData = @{}
for ($r = 2; $r -le 9; $r++) {
    $oldData[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

$mapping = @{ 2 = 3; 3 = 6; 4 = 9; 5 = 2; 6 = 7; 7 = 5; 8 = 4; 9 = 8 }

foreach ($newRow in 2..9) {
    $srcRow = $mapping[$newRow]
    $src = $oldData[$srcRow]

    $ws.Cells.Item($newRow, 4).Value2 = $src.D
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 14).Value2 = $src.N
    $ws.Cells.Item($newRow, 15).Value2 = $src.O
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
    $ws.Cells.Item($newRow, 17).Value2 = $src.Q
    $ws.Cells.Item($newRow, 19).Value2 = $src.S
}
